$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the two top-level picture shapes that were added previously
# (the "powerpoint-icon.jpg" picture and the "ppt_logo.png" picture).
# Because deleting a shape shifts the later items down by one position,
# repeatedly removing item 1 removes both pictures, which sit at the
# very front of the shape collection (ahead of "Group 24").
$s.Shapes.Item(1).Delete()
$s.Shapes.Item(1).Delete()

# The remaining shape is "Group 24"; move it up so its vertical offset
# matches off y="1815259" EMU (142.9338 points == 1815259 EMU once
# PowerPoint rounds the single-precision Top value back to EMU).
$grp = $s.Shapes.Item(1)
$grp.Top = 142.9338
